# Applies the "Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta" edit:
#  - Updates the account-summary header numbers (total mora, worker count, period count)
#  - Replaces the worker detail table with the new set of workers/periods
#  - Removes the now-unused trailing rows, shifting the signature rows up

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values -------------------------------------------------
$ws.Range("E11").Value = 128621   # VALOR MORA
$ws.Range("C13").Value = 3        # Cant. Trabajadores
$ws.Range("F13").Value = 4        # Cant. Periodos

# --- Remove the rows for the periods/workers that are no longer reported ---
# (old rows 19-25 go away; the former "total" row 26 - with its distinct
#  border styling - shifts up into row 19, keeping its own formatting)
$ws.Rows("19:25").Delete()

# --- Rewrite the remaining detail rows with the new worker data ------------
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1047433905"
$ws.Range("D16").Value = "MARY LUZ NAVARRO JIMENEZ"
$ws.Range("E16").Value = "1607"
$ws.Range("F16").Value = 32000
$ws.Range("G16").Value = 737717

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1047433905"
$ws.Range("D17").Value = "MARY LUZ NAVARRO JIMENEZ"
$ws.Range("E17").Value = "1608"
$ws.Range("F17").Value = 32000
$ws.Range("G17").Value = 737717

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1051815149"
$ws.Range("D18").Value = "XAVIER LUIS MENDOZA BUELVAS"
$ws.Range("E18").Value = "1712"
$ws.Range("F18").Value = 29509
$ws.Range("G18").Value = 908526

# Row 19 already holds "VICTOR HUGO TORRES MUNIZ" / periodo 2001 (it was the
# former total row, now shifted up) - only its mora amount changes.
$ws.Range("G19").Value = 877803
